$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 525.4545000000001
$ws.Range("I58").Value = 147.5
$ws.Range("J58").Value = 1533.3334
$ws.Range("K58").Value = 442.5
$ws.Range("L58").Value = 4600.0002
$ws.Range("M58").Value = -292.5
$ws.Range("N58").Value = -4900.0002
$ws.Range("H63").Value = 40000
$ws.Range("J63").Value = 40000
$ws.Range("L63").Value = 40000
$ws.Range("N63").Value = -41248
$ws.Range("H64").Value = 3415.3845
$ws.Range("I64").Value = 3278.5715
$ws.Range("J64").Value = 3575
$ws.Range("K64").Value = 3278.5715
$ws.Range("L64").Value = 3575
$ws.Range("M64").Value = -3030.5715
$ws.Range("N64").Value = -4071
$ws.Range("H66").Value = 40000
$ws.Range("J66").Value = 40000
$ws.Range("L66").Value = 120000
$ws.Range("N66").Value = -126240
$ws.Range("H67").Value = 3415.3845
$ws.Range("I67").Value = 3278.5715
$ws.Range("J67").Value = 3575
$ws.Range("K67").Value = 3278.5715
$ws.Range("L67").Value = 3575
$ws.Range("M67").Value = -2420.5715
$ws.Range("N67").Value = -5291

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3379.7778
$ws.Range("I63").Value = 2365
$ws.Range("J63").Value = 4648.25
$ws.Range("K63").Value = 2365
$ws.Range("L63").Value = 4648.25
$ws.Range("M63").Value = -1679
$ws.Range("N63").Value = -6020.25
$ws.Range("H66").Value = 3379.7778
$ws.Range("I66").Value = 2365
$ws.Range("J66").Value = 4648.25
$ws.Range("K66").Value = 11825
$ws.Range("L66").Value = 23241.25
$ws.Range("M66").Value = -8393
$ws.Range("N66").Value = -30105.25
$ws.Range("H132").Value = 22729938
$ws.Range("I132").Value = 35715830
$ws.Range("K132").Value = 107147490
$ws.Range("M132").Value = -107144960

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2214.2
$ws.Range("I86").Value = 2119.2856
$ws.Range("K86").Value = 2119.2856
$ws.Range("M86").Value = -996.2856000000002
$ws.Range("H89").Value = 2214.2
$ws.Range("I89").Value = 2119.2856
$ws.Range("K89").Value = 10596.428
$ws.Range("M89").Value = -4980.428

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5091.6616
$ws.Range("I31").Value = 2440.5454
$ws.Range("J31").Value = 5631.7036
$ws.Range("K31").Value = 2440.5454
$ws.Range("L31").Value = 5631.7036
$ws.Range("M31").Value = -2145.5454
$ws.Range("N31").Value = -6221.7036
$ws.Range("H34").Value = 5091.6616
$ws.Range("I34").Value = 2440.5454
$ws.Range("J34").Value = 5631.7036
$ws.Range("K34").Value = 2440.5454
$ws.Range("L34").Value = 5631.7036
$ws.Range("M34").Value = -2238.5454
$ws.Range("N34").Value = -6035.7036
$ws.Range("H140").Value = 20919
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 20919
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 20919
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = -31279

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 4200
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 4200
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 12600
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -14098
$ws.Range("H64").Value = 2138
$ws.Range("I64").Value = 1400
$ws.Range("J64").Value = 2507
$ws.Range("K64").Value = 4200
$ws.Range("L64").Value = 7521
$ws.Range("M64").Value = -3930
$ws.Range("N64").Value = -8061
$ws.Range("H66").Value = 4200
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 4200
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 37800
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -45288
$ws.Range("H67").Value = 2138
$ws.Range("I67").Value = 1400
$ws.Range("J67").Value = 2507
$ws.Range("K67").Value = 4200
$ws.Range("L67").Value = 7521
$ws.Range("M67").Value = -3264
$ws.Range("N67").Value = -9393
$ws.Range("H68").Value = 5833.8335
$ws.Range("I68").Value = 1000
$ws.Range("J68").Value = 6800.6
$ws.Range("K68").Value = 3000
$ws.Range("L68").Value = 20401.8
$ws.Range("M68").Value = -2189
$ws.Range("N68").Value = -22023.8
$ws.Range("H69").Value = 73531656
$ws.Range("J69").Value = 73531656
$ws.Range("L69").Value = 220594968
$ws.Range("N69").Value = -220596590
$ws.Range("H70").Value = 5477.4614
$ws.Range("I70").Value = 3006
$ws.Range("J70").Value = 5926.8184
$ws.Range("K70").Value = 9018
$ws.Range("L70").Value = 17780.4552
$ws.Range("M70").Value = -8703
$ws.Range("N70").Value = -18410.4552
$ws.Range("H71").Value = 5833.8335
$ws.Range("I71").Value = 1000
$ws.Range("J71").Value = 6800.6
$ws.Range("K71").Value = 9000
$ws.Range("L71").Value = 61205.4
$ws.Range("M71").Value = -4944
$ws.Range("N71").Value = -69317.39999999999
$ws.Range("H72").Value = 73531656
$ws.Range("J72").Value = 73531656
$ws.Range("L72").Value = 661784904
$ws.Range("N72").Value = -661793016
$ws.Range("H73").Value = 5477.4614
$ws.Range("I73").Value = 3006
$ws.Range("J73").Value = 5926.8184
$ws.Range("K73").Value = 9018
$ws.Range("L73").Value = 17780.4552
$ws.Range("M73").Value = -7926
$ws.Range("N73").Value = -19964.4552
$ws.Range("H81").Value = 1779.3
$ws.Range("I81").Value = 965.5
$ws.Range("J81").Value = 3000
$ws.Range("K81").Value = 2896.5
$ws.Range("L81").Value = 9000
$ws.Range("M81").Value = -1773.5
$ws.Range("N81").Value = -11246
$ws.Range("H84").Value = 1779.3
$ws.Range("I84").Value = 965.5
$ws.Range("J84").Value = 3000
$ws.Range("K84").Value = 8689.5
$ws.Range("L84").Value = 27000
$ws.Range("M84").Value = -3073.5
$ws.Range("N84").Value = -38232
$ws.Range("H124").Value = 1720.25
$ws.Range("J124").Value = 1751.5385
$ws.Range("L124").Value = 5254.6155
$ws.Range("N124").Value = -15074.6155

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 10001420
$ws.Range("I107").Value = 780
$ws.Range("K107").Value = 2340
$ws.Range("M107").Value = -420
